# Updates to .net 5
$p = $ppt.ActivePresentation

# --- 1. Notes master "footer" date placeholder: re-cache the auto date field text ---
$nm = $p.NotesMaster
$dt = $nm.HeadersFooters.DateAndTime
$dt.Text = "12/15/2020"

# --- 2. Slide 11 ("Demo"): content placeholder text updates ---
$s = $p.Slides.Item(11)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

# Paragraph 1: "Activities for W3C TraceContext propagation"
#           -> "Activities for W3C propagation"
$para1 = $tr.Paragraphs(1, 1)
$sub1 = $para1.Characters(16, 28)
$sub1.Text = "W3C propagation"

# Paragraph 3: "OpenTelemetry for Activity adapter"
#           -> "OpenTelemetry for ActivityListener"
$para3 = $tr.Paragraphs(3, 1)
$sub3 = $para3.Characters(19, 16)
$sub3.Text = "ActivityListener"
